# Add a new "Project MM19" milestone block (rows 38-44) to the worksheet,
# mirroring the existing "Project MM18" block, plus a new date number
# format (DD/MM/YY) used for the two new date rows, and move the frozen
# pane / selection down to the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: section header "Project MM19" across A:F, with "Start of Project"
# repeated B:F, style matches existing header rows (style index 1 in xlsx,
# i.e. plain cell style with alignment applied).
$ws.Range("A38").Value = "Project MM19"
$ws.Range("B38").Value = "Start of Project"
$ws.Range("C38").Value = "Start of Project"
$ws.Range("D38").Value = "Start of Project"
$ws.Range("E38").Value = "Start of Project"
$ws.Range("F38").Value = "Start of Project"

# Row 39: "Project MM19 CP" with sample text entries
$ws.Range("A39").Value = "Project MM19 CP"
$ws.Range("B39").Value = "wood"
$ws.Range("C39").Value = "nail"
$ws.Range("D39").Value = "wood"
$ws.Range("E39").Value = "hammer"
$ws.Range("F39").Value = "screw driver"

# Row 40: "Project MM19 Original Baseline" with date values, formatted DD/MM/YY
$ws.Range("A40").Value = "Project MM19 Original Baseline"
$ws.Range("B40").Value = [DateTime]"2000-01-01"
$ws.Range("C40").Value = [DateTime]"2000-01-02"
$ws.Range("D40").Value = [DateTime]"2000-01-03"
$ws.Range("E40").Value = [DateTime]"2000-01-04"
$ws.Range("F40").Value = [DateTime]"2000-01-05"
$ws.Range("B40:G41").NumberFormat = "DD/MM/YY"

# Row 41: "Project MM19 Forecast - Actual" with the same date values
$ws.Range("A41").Value = "Project MM19 Forecast - Actual"
$ws.Range("B41").Value = [DateTime]"2000-01-01"
$ws.Range("C41").Value = [DateTime]"2000-01-02"
$ws.Range("D41").Value = [DateTime]"2000-01-03"
$ws.Range("E41").Value = [DateTime]"2000-01-04"
$ws.Range("F41").Value = [DateTime]"2000-01-05"

# Row 42: "Project MM19 Variance" (values left blank)
$ws.Range("A42").Value = "Project MM19 Variance"

# Row 43: "Project MM19 Status" (values left blank)
$ws.Range("A43").Value = "Project MM19 Status"

# Row 44: "Project MM19 Notes" (values left blank)
$ws.Range("A44").Value = "Project MM19 Notes"

# Touch G38:G41 so the new column G exists in the sheet (empty cells,
# matching the style of the row they are in).
$ws.Range("G38").Value = $ws.Range("G38").Value
$ws.Range("G39").Value = $ws.Range("G39").Value
$ws.Range("G40").Value = $ws.Range("G40").Value
$ws.Range("G41").Value = $ws.Range("G41").Value

# Move the frozen pane / view down to the new block.
$ws.Range("B16").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I40").Select()
